$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 69

$ws.Cells.Item($row, 1).Value = "2025-04-29 09:56:24"
$ws.Cells.Item($row, 2).Value = 208
